$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.393568754196167
$ws.Range("B1").Value = 2.175681114196777
$ws.Range("C1").Value = 4.890871047973633
$ws.Range("D1").Value = 3.541317462921143
$ws.Range("E1").Value = 1.219420075416565
